$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65

# Text-like columns: use a quote-prefix so Excel stores them as literal text
# instead of inferring a date/time/number, then reset the style back to
# "Normal" so no stray quotePrefix style sticks to the cell (matches the
# plain, unstyled text cells used elsewhere in this column).
$ws.Range("A$row").Value = "'2025-02-12"
$ws.Range("A$row").Style = "Normal"

$ws.Range("B$row").Value = "'09:16:47"
$ws.Range("B$row").Style = "Normal"

$ws.Range("C$row").Value = "Wednesday"

$ws.Range("D$row").Value = "'06"
$ws.Range("D$row").Style = "Normal"

$ws.Range("E$row").Value = 127917
$ws.Range("F$row").Value = 142033
$ws.Range("G$row").Value = 169216
$ws.Range("H$row").Value = 158792
$ws.Range("I$row").Value = -1
$ws.Range("J$row").Value = 144426
$ws.Range("K$row").Value = -1
$ws.Range("L$row").Value = -1
$ws.Range("M$row").Value = 191641
$ws.Range("N$row").Value = 115030
$ws.Range("O$row").Value = 44937
$ws.Range("P$row").Value = 28528
$ws.Range("Q$row").Value = 65052
$ws.Range("R$row").Value = -1
$ws.Range("S$row").Value = 43192
$ws.Range("T$row").Value = -1
